$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 99.883077
$ws.Range("H2").Value = 299.649231
$ws.Range("I2").Value = 0.3917580942718377
$ws.Range("J2").Value = 0.3917580942718377
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.08181366666667
$ws.Range("N2").Value = 63.245441
$ws.Range("O2").Value = 0.0571606014598545
$ws.Range("P2").Value = 0.0571606014598545
$ws.Range("Q2").Value = 2105.716417767319
$ws.Range("R2").Value = 18951.44775990587
$ws.Range("S2").Value = 0.02239312829534462
$ws.Range("T2").Value = 0.02239312829534463

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 99.883077
$ws.Range("H3").Value = 299.649231
$ws.Range("I3").Value = 0.3917580942718377
$ws.Range("J3").Value = 0.3917580942718377
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 301.6001486666667
$ws.Range("N3").Value = 904.800446
$ws.Range("O3").Value = 0.8177496571571792
$ws.Range("P3").Value = 0.8177496571571792
$ws.Range("Q3").Value = 30124.75087248411
$ws.Range("R3").Value = 271122.757852357
$ws.Range("S3").Value = 0.3203600472793451
$ws.Range("T3").Value = 0.3203600472793452

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 99.883077
$ws.Range("H4").Value = 299.649231
$ws.Range("I4").Value = 0.3917580942718377
$ws.Range("J4").Value = 0.3917580942718377
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 46.13524966666667
$ws.Range("N4").Value = 138.405749
$ws.Range("O4").Value = 0.1250897413829664
$ws.Range("P4").Value = 0.1250897413829664
$ws.Range("Q4").Value = 4608.130694869892
$ws.Range("R4").Value = 41473.17625382902
$ws.Range("S4").Value = 0.04900491869714792
$ws.Range("T4").Value = 0.04900491869714793

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 124.0161413333333
$ws.Range("H5").Value = 372.048424
$ws.Range("I5").Value = 0.4864119993789693
$ws.Range("J5").Value = 0.4864119993789694
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.08181366666667
$ws.Range("N5").Value = 63.245441
$ws.Range("O5").Value = 0.0571606014598545
$ws.Range("P5").Value = 0.0571606014598545
$ws.Range("Q5").Value = 2614.485183248331
$ws.Range("R5").Value = 23530.36664923498
$ws.Range("S5").Value = 0.02780360244179226
$ws.Range("T5").Value = 0.02780360244179226

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 124.0161413333333
$ws.Range("H6").Value = 372.048424
$ws.Range("I6").Value = 0.4864119993789693
$ws.Range("J6").Value = 0.4864119993789694
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 301.6001486666667
$ws.Range("N6").Value = 904.800446
$ws.Range("O6").Value = 0.8177496571571792
$ws.Range("P6").Value = 0.8177496571571792
$ws.Range("Q6").Value = 37403.28666319968
$ws.Range("R6").Value = 336629.5799687971
$ws.Range("S6").Value = 0.3977632457292902
$ws.Range("T6").Value = 0.3977632457292903

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 124.0161413333333
$ws.Range("H7").Value = 372.048424
$ws.Range("I7").Value = 0.4864119993789693
$ws.Range("J7").Value = 0.4864119993789694
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.13524966666667
$ws.Range("N7").Value = 138.405749
$ws.Range("O7").Value = 0.1250897413829664
$ws.Range("P7").Value = 0.1250897413829664
$ws.Range("Q7").Value = 5721.515643109954
$ws.Range("R7").Value = 51493.64078798959
$ws.Range("S7").Value = 0.06084515120788686
$ws.Range("T7").Value = 0.06084515120788687

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 31.06188766666667
$ws.Range("H8").Value = 93.18566300000001
$ws.Range("I8").Value = 0.121829906349193
$ws.Range("J8").Value = 0.121829906349193
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.08181366666667
$ws.Range("N8").Value = 63.245441
$ws.Range("O8").Value = 0.0571606014598545
$ws.Range("P8").Value = 0.0571606014598545
$ws.Range("Q8").Value = 654.8409279235981
$ws.Range("R8").Value = 5893.568351312383
$ws.Range("S8").Value = 0.006963870722717617
$ws.Range("T8").Value = 0.006963870722717619

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 31.06188766666667
$ws.Range("H9").Value = 93.18566300000001
$ws.Range("I9").Value = 0.121829906349193
$ws.Range("J9").Value = 0.121829906349193
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 301.6001486666667
$ws.Range("N9").Value = 904.800446
$ws.Range("O9").Value = 0.8177496571571792
$ws.Range("P9").Value = 0.8177496571571792
$ws.Range("Q9").Value = 9368.269938133966
$ws.Range("R9").Value = 84314.42944320571
$ws.Range("S9").Value = 0.09962636414854381
$ws.Range("T9").Value = 0.09962636414854383

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 31.06188766666667
$ws.Range("H10").Value = 93.18566300000001
$ws.Range("I10").Value = 0.121829906349193
$ws.Range("J10").Value = 0.121829906349193
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.13524966666667
$ws.Range("N10").Value = 138.405749
$ws.Range("O10").Value = 0.1250897413829664
$ws.Range("P10").Value = 0.1250897413829664
$ws.Range("Q10").Value = 1433.047942619621
$ws.Range("R10").Value = 12897.43148357659
$ws.Range("S10").Value = 0.01523967147793156
$ws.Range("T10").Value = 0.01523967147793156
